$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Dr. G.Chandra Sekhar"
$ws.Range("A2").Value = "Dr. M. S. Sushma Susik"
$ws.Range("A3").Value = "Dr. M.Manjula"
$ws.Range("A4").Value = "Dr. S. Harinath Reddy"
$ws.Range("A5").Value = "Dr. Surendra Kumar Alluri"
